# Weekly update: prepend a new week's worth of data (2 rows) for
# "Vega Monumental Concepción" / Zapallo / Camote, pushing all existing
# rows from 379 down to 381 (sheet grows from A1:R402 to A1:R404).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 379-380; everything currently at/after
# row 379 shifts down by two rows.
$ws.Range("A379:A380").EntireRow.Insert()

# New row 379: Camote, "1a (guarda)"
$ws.Range("A379").Value = 11
$ws.Range("B379").Value = "Vega Monumental Concepción"
$ws.Range("C379").Value = "Bíobío"
$ws.Range("D379").Value = 45106
$ws.Range("E379").Value = 8
$ws.Range("F379").Value = 100112045
$ws.Range("G379").Value = "Zapallo"
$ws.Range("H379").Value = "Camote"
$ws.Range("I379").Value = "1a (guarda)"
$ws.Range("J379").Value = 400
$ws.Range("K379").Value = 300
$ws.Range("L379").Value = 350
$ws.Range("M379").Value = 325
$ws.Range("N379").Value = "$/kilo (volumen en unidades)"
$ws.Range("O379").Value = "Región de O'Higgins"
$ws.Range("P379").Value = 325
$ws.Range("Q379").Value = 1
$ws.Range("R379").Value = "Hortaliza"

# New row 380: Camote, "2a (guarda)"
$ws.Range("A380").Value = 11
$ws.Range("B380").Value = "Vega Monumental Concepción"
$ws.Range("C380").Value = "Bíobío"
$ws.Range("D380").Value = 45106
$ws.Range("E380").Value = 8
$ws.Range("F380").Value = 100112045
$ws.Range("G380").Value = "Zapallo"
$ws.Range("H380").Value = "Camote"
$ws.Range("I380").Value = "2a (guarda)"
$ws.Range("J380").Value = 200
$ws.Range("K380").Value = 250
$ws.Range("L380").Value = 250
$ws.Range("M380").Value = 250
$ws.Range("N380").Value = "$/kilo (volumen en unidades)"
$ws.Range("O380").Value = "Región de O'Higgins"
$ws.Range("P380").Value = 250
$ws.Range("Q380").Value = 1
$ws.Range("R380").Value = "Hortaliza"
